$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric need the Text number
# format applied first, otherwise Excel COM auto-converts them to real
# numbers instead of leaving them as literal text (matching the source data).
$textForceCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin name / link / price / change values.
$ws.Range("D2").Value = "27.026.17"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "1.818.16"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -1.35%  "
$ws.Range("D5").Value = "310.48"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("D7").Value = "0.4211"
$ws.Range("E7").Value = "  -2.17%  "
$ws.Range("D8").Value = "0.3661"
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").Value = "0.07200"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "0.8413"
$ws.Range("E10").Value = "  -3.73%  "
$ws.Range("D11").Value = "20.80"
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("D12").Value = "1.815.30"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "6.620"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").Value = "0.07065"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "5.266"
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("D16").Value = "88.85"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "0.000008799"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("B20").Value = "BitDAO"
$ws.Range("C20").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D20").Value = "0.5046"
$ws.Range("E20").Value = "  -3.73%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "14.94"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "27.100.12"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.110"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "10.80"
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.046.02"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "1.976"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "151.50"
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.228"
$ws.Range("E28").Value = "  +4.01%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "18.23"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "5.201"
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "115.78"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.08781"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "1.173"
$ws.Range("E33").Value = "  -3.95%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.953"
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7380"
$ws.Range("E35").Value = "  -4.59%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "4.408"
$ws.Range("E36").Value = "  -2.90%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.092"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01955"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.05226"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "7.258"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "2.877"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.1683"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.5027"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "8.559"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.50"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.4730"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "106.08"
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.06359"
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").Value = "1.879"
$ws.Range("E51").Value = "  +2.38%  "
